$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Fitness values for rows 2-87 (C2:C87), as described by the commit
# "correction in sa algorithm and 746 logs".
$values = @(
    8297, 8297, 8297, 8297, 8297, 8297, 8297, 8297, 8297, 8297,
    8297, 8297, 8297, 8297, 8297, 8297, 8297, 8297, 8297, 8297,
    8297, 8297, 8297, 8145, 8145, 8145, 8145, 7979, 7979, 7979,
    7979, 7979, 7979, 7979, 7979, 7979, 7979, 7946, 7946, 7946,
    7946, 7946, 7946, 7657, 7657, 7657, 7657, 7657, 7657, 7657,
    7657, 7657, 7573, 7573, 7573, 7573, 7573, 7573, 7573, 7573,
    7573, 7573, 7573, 7573, 7573, 7573, 7573, 7573, 7573, 7573,
    7573, 7573, 7573, 7573, 7573, 7573, 7573, 7573, 7573, 7573,
    7573, 7573, 7573, 7573, 7573, 7573
)

$startRow = 2
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $values[$i]
}
